$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario Mapping")

# --- D column: fill in scenario-mapping descriptions (order matters for
#     shared-string table layout) ---
$ws.Range("D32").Value = '"contentChanged": false,
"statusName": "Suspended",'
$ws.Range("D33").Value = '"contentChanged": false,
"statusName": "Suspended",
WHEN EXIST addProduct items  AND
WHEN EXIST removeProduct item AND'
$ws.Range("D31").Value = 'Cell1:
"contentChanged": false,
"statusName": "Withdrawn",
Cell2:
"contentChanged": false
"statusName": Suspended 
-------------------------------------
WHEN EXIST addProduct items  AND
WHEN EXIST removeProduct item AND'
$ws.Range("D30").Value = '"contentChanged": true,
"statusName": "Withdrawn",
WHEN EXIST addProduct items  AND
WHEN EXIST removeProduct item AND
CANCEL AVCS UOS WHEN STATUS="NotForSale"'
$ws.Range("D29").Value = '"contentChanged": true,
"statusName": "Suspended",'
$ws.Range("D28").Value = '"contentChange": true,
"statusName": "Update",
'
$ws.Range("D27").Value = '"contentChanged": true,
"isNewCell": true,'
$ws.Range("D26").Value = 'isNewCell = true & IsNewUnitOfSale = true
WHEN StatusName = “New Edition” AND
WHEN EXIST addProduct items  AND
WHEN EXIST removeProduct item AND'
$ws.Range("D25").Value = 'contentChanged = false'
$ws.Range("D24").Value = 'If “contentChanged”:false
WHEN EXIST addProduct items  AND
WHEN EXIST removeProduct item AND'
$ws.Range("D23").Value = 'If “contentChanged”:false
isNewCell=true,
inUnitOfSales contains multiple unit of sales with unitOfSaleType="unit",
WHEN EXIST addProduct items  AND
WHEN EXIST removeProduct item AND'
$ws.Range("D22").Value = 'isNewCell = true  AND
IsNewUnitOfSale = False AND
WHEN EXIST addProduct items  AND
WHEN EXIST replacedBy items AND
WHEN EXIST removeProduct item AND
WHEN StatusName = “Cancellation Update” AND
WHEN status = “NotForSale”'
$ws.Range("D21").Value = '"contentChanged": true,
"statusName": "New Edition",
"isNewCell": false,'
$ws.Range("D20").Value = '"contentChanged": true,
WHEN EXIST addProduct items  AND
WHEN EXIST replacedBy items AND
WHEN EXIST removeProduct item AND'
$ws.Range("D19").Value = '"contentChanged": true,
"statusName": "New Edition",
WHEN EXIST addProduct items  AND
WHEN EXIST replacedBy items AND
WHEN EXIST removeProduct item AND'

# --- enable wrap text so the new multi-line notes render correctly ---
$ws.Range("D21").WrapText = $true
$ws.Range("D22").WrapText = $true
$ws.Range("D23").WrapText = $true
$ws.Range("D24").WrapText = $true
$ws.Range("D25").WrapText = $true
$ws.Range("D26").WrapText = $true
$ws.Range("D27").WrapText = $true
$ws.Range("D28").WrapText = $true
$ws.Range("D29").WrapText = $true
$ws.Range("D30").WrapText = $true
$ws.Range("D31").WrapText = $true
$ws.Range("D32").WrapText = $true
$ws.Range("D33").WrapText = $true

# --- row heights to fit the new wrapped content ---
$ws.Rows.Item(19).RowHeight = 75
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 45
$ws.Rows.Item(22).RowHeight = 105
$ws.Rows.Item(23).RowHeight = 90
$ws.Rows.Item(24).RowHeight = 45
$ws.Rows.Item(26).RowHeight = 60
$ws.Rows.Item(27).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 45
$ws.Rows.Item(29).RowHeight = 30
$ws.Rows.Item(30).RowHeight = 75
$ws.Rows.Item(31).RowHeight = 135
$ws.Rows.Item(32).RowHeight = 30
$ws.Rows.Item(33).RowHeight = 60

# --- last selection before save ---
$ws.Range("D20").Select() | Out-Null

Write-Output "done"